# Add margin and axis label funs and line drawing for pub ready charts
#
# This reorders Sheet30/Sheet31 (moving Sheet30 before Sheet31) and adds a
# brand new Sheet32 at the end containing GDP growth-component data
# (percentage-point contributions), wiring up the selections/tab exactly as
# produced by Excel when a worksheet is added & made active.

$wb = $excel.ActiveWorkbook

# --- 1. Reorder Sheet30 to sit before Sheet31 -----------------------------
$sheet30 = $wb.Worksheets.Item("Sheet30")
$sheet31 = $wb.Worksheets.Item("Sheet31")
$sheet30.Move($sheet31)

# Tidy up the selections left behind on the two existing sheets (re-fetch
# fresh references since their tab position just changed).
$sheet30 = $wb.Worksheets.Item("Sheet30")
$sheet31 = $wb.Worksheets.Item("Sheet31")

$sheet30.Select()
$sheet30.Range("N25").Select()

$sheet31.Select()
$sheet31.Range("A1:S7").Select()
$sheet31.Range("N2").Activate()

# --- 2. Add the new Sheet32 after Sheet31 ---------------------------------
$sheet31 = $wb.Worksheets.Item("Sheet31")
$newSheet = $wb.Worksheets.Add($null, $sheet31)
$newSheet.Name = "Sheet32"

# Header row (same headers as the other data sheets on this tab group).
$newSheet.Range("A1").Value = "serija"
$newSheet.Range("B1").Value = "enota"
$newSheet.Range("C1").Value = "legenda"
$newSheet.Range("D1").Value = "barva"
$newSheet.Range("E1").Value = "tip"
$newSheet.Range("F1").Value = "stacked"
$newSheet.Range("G1").Value = "drseca_obdobja"
$newSheet.Range("H1").Value = "drseca_poravnava"
$newSheet.Range("I1").Value = "rast"
$newSheet.Range("J1").Value = "indeks_obdobje"
$newSheet.Range("K1").Value = "velikost"
$newSheet.Range("L1").Value = "naslov"
$newSheet.Range("M1").Value = "xmin"
$newSheet.Range("N1").Value = "xmax"
$newSheet.Range("O1").Value = "opomba"
$newSheet.Range("P1").Value = "stolpci_legende"
$newSheet.Range("Q1").Value = "datum_podatkov"
$newSheet.Range("R1").Value = "leva_y_os"
$newSheet.Range("S1").Value = "desna_y_os"

# Data rows: quarterly GDP growth broken into percentage-point contributions.
$newSheet.Range("A2").Value = "SURS--0300230S--B1GQ--G4--N--Q"
$newSheet.Range("B2").Value = "odstotne točke"
$newSheet.Range("C2").Value = "Bruto domači proizvod"
$newSheet.Range("E2").Value = "line"
$newSheet.Range("M2").Value = 40909
$newSheet.Range("M2").NumberFormat = "m/d/yyyy"
$newSheet.Range("N2").Value = 45292
$newSheet.Range("N2").NumberFormat = "m/d/yyyy"

$newSheet.Range("A3").Value = "SURS--0300230S--P3_S13--GO4--N--Q"
$newSheet.Range("B3").Value = "odstotne točke"
$newSheet.Range("C3").Value = "....Končna potrošnja države "
$newSheet.Range("E3").Value = "line"

$newSheet.Range("A4").Value = "SURS--0300230S--P31_S14_D--GO4--N--Q"
$newSheet.Range("B4").Value = "odstotne točke"
$newSheet.Range("C4").Value = "....Končna potrošnja gospodinjstev"
$newSheet.Range("E4").Value = "line"

$newSheet.Range("A5").Value = "SURS--0300230S--P51G--GO4--N--Q"
$newSheet.Range("B5").Value = "odstotne točke"
$newSheet.Range("C5").Value = "....Bruto investicije v osnovna sredstva"
$newSheet.Range("E5").Value = "line"

$newSheet.Range("A6").Value = "SURS--0300230S--P52--GO4--N--Q"
$newSheet.Range("B6").Value = "odstotne točke"
$newSheet.Range("C6").Value = "....Spremembe zalog"
$newSheet.Range("E6").Value = "line"

$newSheet.Range("A7").Value = "SURS--0300230S--B11--GO4--N--Q"
$newSheet.Range("B7").Value = "odstotne točke"
$newSheet.Range("C7").Value = "Saldo menjave s tujino "
$newSheet.Range("E7").Value = "line"

# Make Sheet32 the active tab with the same selection Excel leaves behind
# after a column-fill operation.
$newSheet.Select()
$newSheet.Range("D2:D7").Select()
$newSheet.Range("D2").Activate()
